# Merge the cell description function by Prerna.
#
# The "check student's answer for zipcode of New York" description (row for
# identifier Z2 / cell B5) is removed from the Zipcode_CheckOrder sheet - its
# description gets merged/cleared out, leaving the Hidden? column value (D3)
# as the only remaining data for that row. The now-unused shared string is
# pruned automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zipcode_CheckOrder")

# Clear the "Descriptions? (Optional)" cell for the New York row (Z2 / B5).
$ws.Range("C3").ClearContents()

# Excel leaves the selection on the cell that was just edited.
$ws.Range("C3").Select()
